$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 62
$ws.Range("H62").Value = 2229.9524
$ws.Range("I62").Value = 2461.7273
$ws.Range("J62").Value = 1975
$ws.Range("K62").Value = 2461.7273
$ws.Range("L62").Value = 1975
$ws.Range("M62").Value = -1837.7273
$ws.Range("N62").Value = -3223

# ALC row 65
$ws.Range("H65").Value = 2229.9524
$ws.Range("I65").Value = 2461.7273
$ws.Range("J65").Value = 1975
$ws.Range("K65").Value = 12308.6365
$ws.Range("L65").Value = 9875
$ws.Range("M65").Value = -9188.636500000001
$ws.Range("N65").Value = -16115

$ws = $wb.Worksheets.Item("ARM")
# ARM row 61
$ws.Range("H61").Value = 26369746
$ws.Range("I61").Value = 37075096
$ws.Range("K61").Value = 37075096
$ws.Range("M61").Value = -37074884

# ARM row 63
$ws.Range("H63").Value = 5585.357
$ws.Range("I63").Value = 4476.5386
$ws.Range("J63").Value = 20000
$ws.Range("K63").Value = 4476.5386
$ws.Range("L63").Value = 20000
$ws.Range("M63").Value = -3790.5386
$ws.Range("N63").Value = -21372

# ARM row 66
$ws.Range("H66").Value = 5585.357
$ws.Range("I66").Value = 4476.5386
$ws.Range("J66").Value = 20000
$ws.Range("K66").Value = 22382.693
$ws.Range("L66").Value = 100000
$ws.Range("M66").Value = -18950.693
$ws.Range("N66").Value = -106864

# ARM row 110
$ws.Range("H110").Value = 1034.2693
$ws.Range("I110").Value = 705.8570999999999
$ws.Range("J110").Value = 2413.6
$ws.Range("K110").Value = 705.8570999999999
$ws.Range("L110").Value = 2413.6
$ws.Range("M110").Value = 1339.1429
$ws.Range("N110").Value = -6503.6

# ARM row 122
$ws.Range("H122").Value = 3878
$ws.Range("I122").Value = 3012
$ws.Range("J122").Value = 4166.6665
$ws.Range("K122").Value = 9036
$ws.Range("L122").Value = 12499.9995
$ws.Range("M122").Value = -6586
$ws.Range("N122").Value = -17399.9995

# ARM row 136
$ws.Range("H136").Value = 26369746
$ws.Range("I136").Value = 37075096
$ws.Range("K136").Value = 111225288
$ws.Range("M136").Value = -111222738

$ws = $wb.Worksheets.Item("BSM")
# BSM row 55
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# CRP row 99
$ws.Range("H99").Value = 1763.3334
$ws.Range("I99").Value = 1700
$ws.Range("K99").Value = 1700
$ws.Range("M99").Value = -202

# CRP row 126
$ws.Range("H126").Value = 1763.3334
$ws.Range("I126").Value = 1700
$ws.Range("K126").Value = 5100
$ws.Range("M126").Value = -2630

# CRP row 134
$ws.Range("H134").Value = 30912.916
$ws.Range("I134").Value = 746.1053000000001
$ws.Range("K134").Value = 2238.3159
$ws.Range("M134").Value = 296.6840999999999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 88
$ws.Range("H88").Value = 2968.75
$ws.Range("J88").Value = 2968.75
$ws.Range("L88").Value = 8906.25
$ws.Range("N88").Value = -9762.25

# CUL row 91
$ws.Range("H91").Value = 2968.75
$ws.Range("J91").Value = 2968.75
$ws.Range("L91").Value = 8906.25
$ws.Range("N91").Value = -11870.25

# CUL row 137
$ws.Range("H137").Value = 1800.35
$ws.Range("J137").Value = 2561.5454
$ws.Range("L137").Value = 7684.6362
$ws.Range("N137").Value = -17884.6362

$ws = $wb.Worksheets.Item("GSM")
# GSM row 51
$ws.Range("H51").Value = 50166.668
$ws.Range("J51").Value = 50166.668
$ws.Range("L51").Value = 50166.668
$ws.Range("N51").Value = -51184.668

# GSM row 80
$ws.Range("H80").Value = 3233.75
$ws.Range("I80").Value = 1766.6666
$ws.Range("J80").Value = 3572.3076
$ws.Range("K80").Value = 1766.6666
$ws.Range("L80").Value = 3572.3076
$ws.Range("M80").Value = -768.6666
$ws.Range("N80").Value = -5568.3076

# GSM row 83
$ws.Range("H83").Value = 3233.75
$ws.Range("I83").Value = 1766.6666
$ws.Range("J83").Value = 3572.3076
$ws.Range("K83").Value = 8833.333000000001
$ws.Range("L83").Value = 17861.538
$ws.Range("M83").Value = -3841.333000000001
$ws.Range("N83").Value = -27845.538

# GSM row 102
$ws.Range("H102").Value = 1996
$ws.Range("I102").Value = 1892.5714
$ws.Range("J102").Value = 2140.8
$ws.Range("K102").Value = 1892.5714
$ws.Range("L102").Value = 2140.8
$ws.Range("M102").Value = -270.5714
$ws.Range("N102").Value = -5384.8

# GSM row 122
$ws.Range("H122").Value = 2500.5386
$ws.Range("I122").Value = 1845.2222
$ws.Range("J122").Value = 3975
$ws.Range("K122").Value = 5535.6666
$ws.Range("L122").Value = 11925
$ws.Range("M122").Value = -3085.6666
$ws.Range("N122").Value = -16825

# GSM row 126
$ws.Range("H126").Value = 2273.625
$ws.Range("I126").Value = 1666.3334
$ws.Range("J126").Value = 2638
$ws.Range("K126").Value = 4999.0002
$ws.Range("L126").Value = 7914
$ws.Range("M126").Value = -2529.0002
$ws.Range("N126").Value = -12854

# GSM row 132
$ws.Range("H132").Value = 78904.80499999999
$ws.Range("I132").Value = 57217.055
$ws.Range("J132").Value = 127702.25
$ws.Range("K132").Value = 171651.165
$ws.Range("L132").Value = 383106.75
$ws.Range("M132").Value = -169121.165
$ws.Range("N132").Value = -388166.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 2500.5715
$ws.Range("I40").Value = 2541.9167
$ws.Range("J40").Value = 2252.5
$ws.Range("K40").Value = 2541.9167
$ws.Range("L40").Value = 2252.5
$ws.Range("M40").Value = -2405.9167
$ws.Range("N40").Value = -2524.5

# LTW row 122
$ws.Range("H122").Value = 3789.8235
$ws.Range("I122").Value = 4299.6
$ws.Range("J122").Value = 3387.3684
$ws.Range("K122").Value = 12898.8
$ws.Range("L122").Value = 10162.1052
$ws.Range("M122").Value = -10448.8
$ws.Range("N122").Value = -15062.1052

# LTW row 132
$ws.Range("H132").Value = 40984.52
$ws.Range("I132").Value = 2299.775
$ws.Range("J132").Value = 169933.67
$ws.Range("K132").Value = 6899.325000000001
$ws.Range("L132").Value = 509801.01
$ws.Range("M132").Value = -4369.325000000001
$ws.Range("N132").Value = -514861.01

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 1451
$ws.Range("I122").Value = 1402
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4206
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1756
$ws.Range("N122").Value = -9400

# WVR row 126
$ws.Range("H126").Value = 1255.5714
$ws.Range("I126").Value = 1255.5714
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3766.7142
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1296.7142
$ws.Range("N126").ClearContents()
